$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# "Estado de Cuenta" report: a new mora period (2509) is added below the
# existing 2507 / 2508 rows, the totals are refreshed, and the signature
# block shifts down by one row to make room.
# ---------------------------------------------------------------------------

# 1) Insert a new row 18 (pushes the signature block at rows 22-23 down to
#    23-24, same as Excel's native "Insert Copied Cells" workflow).
$ws.Rows.Item(18).Insert()

# 2) Duplicate the old last data row (17 -> 2508, "closing" box-border style)
#    down into the freshly inserted row 18.
$ws.Range("B17:J17").Copy($ws.Range("B18:J18"))

# 3) The old row 17 is no longer the last row, so it now takes on the plain
#    "interior" row look that row 16 already has (copy values+format from
#    row 16, then restore row 17's own "2508" label).
$ws.Range("B16:J16").Copy($ws.Range("B17:J17"))
$ws.Range("E17").Value = "2508"

# 4) The brand-new row 18 represents period 2509.
$ws.Range("E18").Value = "2509"

# 5) Refresh the header summary: one more "Cant. Periodos" and the mora
#    total grows by the new period's value (56940).
$ws.Range("F13").Value = 3
$ws.Range("E11").Value = 170820
